$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the front for "chapter_id" and shift existing
# question/option/answer/etc. columns one place to the right.
$ws.Columns("A:A").Insert()

$ws.Range("A1").Value = "chapter_id"
$ws.Range("J1").Value = "time_duration"

# Fill chapter_id values for data rows 2-30 with the generated uuid.
$uuid = "f1a3a82b-41d3-4d63-bf0d-a30a85538249"
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value = $uuid
}

# Column widths (characters -> matches header text length + 30, converted
# through Excel's internal width grid as closely as this runtime allows).
$ws.Columns("A:A").ColumnWidth = 39.83333333333333
$ws.Columns("B:B").ColumnWidth = 37.83333333333333
$ws.Columns("C:F").ColumnWidth = 36.83333333333333
$ws.Columns("G:G").ColumnWidth = 35.83333333333333
$ws.Columns("H:H").ColumnWidth = 46.83333333333333
$ws.Columns("I:I").ColumnWidth = 34.83333333333333
$ws.Columns("J:J").ColumnWidth = 42.83333333333333
